$d = $word.ActiveDocument

# "remove personal info from website" — strip the email address and phone
# number (plus their " | " separators) from the header/contact line, while
# leaving the LinkedIn and GitHub-pages hyperlinks that follow untouched.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("joshuaalev@gmail.com | 818.257.4496 | ", $true, $false, $false, `
              $false, $false, $true, 1, $false, "", 2)
